$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 219 (pushes existing rows 219-308 down to 221-310)
$ws.Rows.Item(219).Insert()
$ws.Rows.Item(219).Insert()

# Populate new row 219
$ws.Range('A219').Value = 5
$ws.Range('B219').Value = 'Macroferia Regional de Talca'
$ws.Range('C219').Value = 'Maule'
$ws.Range('D219').Value = 44795
$ws.Range('E219').Value = 7
$ws.Range('F219').Value = 100112008
$ws.Range('G219').Value = 'Coliflor'
$ws.Range('H219').Value = 'Sin especificar'
$ws.Range('I219').Value = 'Primera'
$ws.Range('J219').Value = 2000
$ws.Range('K219').Value = 1000
$ws.Range('L219').Value = 1000
$ws.Range('M219').Value = 1000
$ws.Range('N219').Value = '$/unidad'
$ws.Range('O219').Value = 'Región del Maule'
$ws.Range('P219').Value = 1000
$ws.Range('Q219').Value = 1
$ws.Range('R219').Value = 'Hortaliza'

# Populate new row 220
$ws.Range('A220').Value = 5
$ws.Range('B220').Value = 'Macroferia Regional de Talca'
$ws.Range('C220').Value = 'Maule'
$ws.Range('D220').Value = 44795
$ws.Range('E220').Value = 7
$ws.Range('F220').Value = 100112008
$ws.Range('G220').Value = 'Coliflor'
$ws.Range('H220').Value = 'Sin especificar'
$ws.Range('I220').Value = 'Segunda'
$ws.Range('J220').Value = 2000
$ws.Range('K220').Value = 800
$ws.Range('L220').Value = 800
$ws.Range('M220').Value = 800
$ws.Range('N220').Value = '$/unidad'
$ws.Range('O220').Value = 'Región del Maule'
$ws.Range('P220').Value = 800
$ws.Range('Q220').Value = 1
$ws.Range('R220').Value = 'Hortaliza'
